# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet (positioned between "2021-Q4" and
#    "总计") with the quarter's per-fund holding detail.
# 2. Insert a new leading row into "总计" summarising the 2022-Q1 totals,
#    pushing the existing 2021-Q4 / 2021-Q3 rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q1" sheet (Worksheets.Add() inserts at the
#    front; populate it there, then Move() it into place right before
#    "总计" so the final tab order is 2021-Q3, 2021-Q4, 2022-Q1, 总计).
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add()
$q1.Name = "2022-Q1"

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'000977"
$q1.Range("C2").Value = "长城环保主题灵活配置混合"
$q1.Range("D2").Value = "'9.63"
$q1.Range("E2").Value = "'81.27"
$q1.Range("F2").Value = "'2.69"
$q1.Range("G2").Value = "'0.2590"
$q1.Range("H2").Value = 10

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'010049"
$q1.Range("C3").Value = "长城成长先锋混合A"
$q1.Range("D3").Value = "'7.27"
$q1.Range("E3").Value = "'81.13"
$q1.Range("F3").Value = "'2.65"
$q1.Range("G3").Value = "'0.1927"
$q1.Range("H3").Value = 10

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "'002542"
$q1.Range("C4").Value = "长城久鼎灵活配置混合"
$q1.Range("D4").Value = "'3.81"
$q1.Range("E4").Value = "'81.60"
$q1.Range("F4").Value = "'2.71"
$q1.Range("G4").Value = "'0.1033"
$q1.Range("H4").Value = 10

$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "'010050"
$q1.Range("C5").Value = "长城成长先锋混合C"
$q1.Range("D5").Value = "'0.90"
$q1.Range("E5").Value = "'81.13"
$q1.Range("F5").Value = "'2.65"
$q1.Range("G5").Value = "'0.0238"
$q1.Range("H5").Value = 10

# Match the bold/centred/bordered header + index-column styling used by
# the sibling quarter sheets (copy format from "2021-Q4", then restore
# the values this sheet actually owns since Copy also brings text along).
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Range("B1:H1").Copy($q1.Range("B1:H1"))
$q4.Range("A2:A4").Copy($q1.Range("A2:A5"))

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"
$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1
$q1.Range("A4").Value = 2
$q1.Range("A5").Value = 3

$q1.Move($wb.Worksheets.Item("总计"))

# ---------------------------------------------------------------------
# 2. Insert the 2022-Q1 summary row at the top of "总计", pushing the
#    2021-Q4 / 2021-Q3 rows down to rows 3 / 4.
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

$oldB2 = $tot.Range("B2").Value()
$oldC2 = $tot.Range("C2").Value()
$oldD2 = $tot.Range("D2").Value()
$oldB3 = $tot.Range("B3").Value()
$oldC3 = $tot.Range("C3").Value()
$oldD3 = $tot.Range("D3").Value()

$tot.Range("A3").Copy($tot.Range("A4"))
$tot.Range("A4").Value = 2
$tot.Range("B4").Value = $oldB3
$tot.Range("C4").Value = $oldC3
$tot.Range("D4").Value = $oldD3

$tot.Range("A3").Value = 1
$tot.Range("B3").Value = $oldB2
$tot.Range("C3").Value = $oldC2
$tot.Range("D3").Value = $oldD2

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 4
$tot.Range("D2").Value = 0.58
